# Update the "Correspond Handoff/Handback DateTime" cells for the
# "84699ea6-f550-44bc-814d-c81d41b09f2c" handback entry, for both the
# zh-cn and de-de report sheets, to reflect the freshly re-generated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$sheetNames = @("zh-cn", "de-de")

$handoffTimes = @{
    "zh-cn" = "2016-03-17 16:16:07"
    "de-de" = "2016-03-17 16:16:11"
}

$handbackTimes = @{
    "zh-cn" = "2016-03-17 16:16:35"
    "de-de" = "2016-03-17 16:16:45"
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Rows 3 and 5 both correspond to the 84699ea6-... file in this sheet.
    $ws.Range("E3").Value = $handoffTimes[$name]
    $ws.Range("H3").Value = $handbackTimes[$name]

    $ws.Range("E5").Value = $handoffTimes[$name]
    $ws.Range("H5").Value = $handbackTimes[$name]
}
